# "time and sales updates"
# Insert a new bond record into the "List" sheet, between the existing
# 01/04/23-maturity row (row 63) and the 02/11/23-maturity row (old row 64),
# keeping the sheet's ascending sort by maturity_date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")

# Insert a new row at 64 - Excel shifts rows 64:123 down to 65:124 and
# copies the number formats from the row above (row 63), same as a manual
# "Insert Sheet Rows" in the UI.
$ws.Rows(64).Insert()

# series_mosb - stored as literal text (leading apostrophe forces the
# quote-prefixed text style instead of the date auto-conversion Excel
# would otherwise apply to "01/11/23").
$ws.Range("A64").Value = "'01/11/23"

# isin
$ws.Range("C64").Value = "PIBL1222A022"

# coupon_rate - these placeholder/repo-bill rows use "-" for coupon
$ws.Range("D64").Value = "-"

# issue_date / maturity_date (serial dates, same style as surrounding rows)
$ws.Range("E64").Value = 44573
$ws.Range("F64").Value = 44937

# watchlist - same as the other freshly-entered rows (row 2 & 3), driven by
# a literal FALSE() formula rather than a bare boolean literal
$ws.Range("G64").Formula = "=FALSE()"

# Column C ("isin") width was nudged via AutoFit/bestfit after the new,
# slightly different entry was added.
$ws.Columns("C").ColumnWidth = 12.49

# Selection moved down to where the new row was typed in.
$ws.Range("C123").Select()
